# Remove leftover HTML markup (<small>, <em>) from the English gloss
# column (A) that was never meant to reach the published spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "to go (destination に/へ)"
    3  = "to go back; to return (destination に/へ)"
    4  = "to listen; to hear (～を)"
    5  = "to drink (～を)"
    6  = "to speak; to talk (language を/で)"
    7  = "to read (～を)"
    52 = "to eat (～を)"
    54 = "to see; to look at; to watch (～を)"
    55 = "to come (destination に/へ)"
    56 = "to do (～を)"
    57 = "to study (～を)"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
